# Update the staramr run metadata recorded on the "Settings" sheet.
# This reflects a re-run of the pipeline (new Galaxy job / output paths,
# new start/end timestamps, and a corrected total runtime) after fixing
# a problem with the database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# command_line (B2): new Galaxy job id / dataset output paths
$ws.Range("B2").Value = "/shared/ifbstor1/galaxy/mutable-data/dependencies/_conda/envs/mulled-v1-50d167472d35f4b3d40c0e43369adafb301df8a06c0471a38dc9a39220c30ff7/bin/staramr search --nprocs 1 --genome-size-lower-bound 4000000 --genome-size-upper-bound 6000000 --minimum-N50-value 10000 --minimum-contig-length 300 --unacceptable-number-contigs 1000 --pid-threshold 98.0 --percent-length-overlap-resfinder 60.0 --percent-length-overlap-plasmidfinder 60.0 --percent-length-overlap-pointfinder 95.0 --output-summary /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_93d66200-ea55-4dff-a489-8c02d9212aa7.dat --output-detailed-summary /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_df1ae156-84f5-453d-b381-acaf3ffbcf69.dat --output-resfinder /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_45578c87-faf5-49d4-a841-e6e71f81296c.dat --output-plasmidfinder /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_3049e167-768d-4d80-bba5-9aeb1c569a8b.dat --output-settings /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_acb5634f-be75-4cc0-9204-7728ecc14a4e.dat --output-excel results.xlsx --output-mlst /shared/ifbstor1/galaxy/jobs/003/856/3856697/outputs/galaxy_dataset_92c71ffa-2402-4439-950c-ec3b2d98f435.dat --output-hits-dir staramr_hits shovill_contigs_fasta.fasta"

# start_time (B4)
$ws.Range("B4").Value = "2024-01-25 12:10:23"

# end_time (B5)
$ws.Range("B5").Value = "2024-01-25 12:14:49"

# total_minutes (B6) - keep this a text value (as it was before: "0.43"),
# not a numeric one, by forcing text entry via a leading apostrophe.
$ws.Range("B6").Value = "'4.43"
